$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - first sheet
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1082
$ws1.Range("F3").Value = 101
$ws1.Range("F4").Value = 1657
$ws1.Range("G4").Value = 60
$ws1.Range("F5").Value = 752
$ws1.Range("F6").Value = 162

# Sheet "全部类型" (All types) - fourth sheet
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1082
$ws4.Range("F3").Value = 101
$ws4.Range("F4").Value = 1657
$ws4.Range("G4").Value = 60
$ws4.Range("F6").Value = 752
$ws4.Range("F7").Value = 162
